$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (pushes existing rows 9-12 down to 10-13),
# carrying down formatting the same way Excel's own Insert does.
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the new weekly entry.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44629
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101007
$ws.Range("J9").Value = "Kiwi"
$ws.Range("K9").Value = "Hayward"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 17500
$ws.Range("Q9").Value = "$/bandeja 18 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 972
$ws.Range("T9").Value = 18
